$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header text changes (columns U..AC get new labels) ---
$ws.Range("U1").Value = "Size"
$ws.Range("V1").Value = "Color"
$ws.Range("W1").Value = "Writer Type"
$ws.Range("X1").Value = "Form Factor"
$ws.Range("Y1").Value = "Model Name"
$ws.Range("Z1").Value = "Model ID"
$ws.Range("AA1").Value = "Interface"
$ws.Range("AB1").Value = "Disclaimer"
$ws.Range("AC1").Value = "In Sales Package"

# --- Row 2 value changes ---
$ws.Range("A2").NumberFormat = "General"
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "HP External DVD Writer"
$ws.Range("C2").Value = 2899
$ws.Range("D2").Value = 2000
$ws.Range("E2").Value = 100
$ws.Range("F2").ClearContents()
$ws.Range("G2").Value = "Suitable For All Brand Laptops And Desktops"
$ws.Range("I2").ClearContents()
$ws.Range("K2").ClearContents()
$ws.Range("M2").ClearContents()
$ws.Range("O2").ClearContents()
$ws.Range("Q2").ClearContents()
$ws.Range("R2").ClearContents()
$ws.Range("S2").ClearContents()
$ws.Range("T2").Value = "HP External DVD Writer"
$ws.Range("U2").ClearContents()
$ws.Range("V2").Value = "Black"
$ws.Range("W2").Value = "CD/DVD"
$ws.Range("X2").Value = "Portable"
$ws.Range("Y2").ClearContents()
$ws.Range("Z2").ClearContents()
$ws.Range("AA2").Value = "USB"
$ws.Range("AB2").ClearContents()
$ws.Range("AC2").Value = "1 External DVD Writer"
$ws.Range("AD2").ClearContents()
$ws.Range("AE2").ClearContents()
$ws.Range("AF2").ClearContents()
$ws.Range("AG2").ClearContents()
$ws.Range("AH2").ClearContents()
$ws.Range("AI2").ClearContents()
$ws.Range("AJ2").ClearContents()
$ws.Range("AK2").ClearContents()

# --- View / selection changes ---
$ws.Range("AB1").Select()
